$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$lq = [char]0x2018
$rq = [char]0x2019

# ---------------------------------------------------------------------------
# Helper: find the paragraph whose text starts with $needle and return its
# 1-based paragraph index.
# ---------------------------------------------------------------------------
function Find-ParaIndex($needle) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $ok = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "NOT FOUND: $needle"
        return -1
    }
    $target = $rng.Start
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Start -eq $target) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) Insert a brand-new paragraph before "Do this lab in groups of two..."
#    that references the example document.
# ---------------------------------------------------------------------------
$doThisIdx = Find-ParaIndex("Do this lab in groups of two")
$doThisPara = $d.Paragraphs($doThisIdx)
$doThisPara.Range.InsertParagraphBefore()

$exampleIdx = $doThisIdx
$examplePara = $d.Paragraphs($exampleIdx)
$exampleXml = "<w:p $wNs><w:r><w:t>There is an example of what you need to do for this lab in $lq</w:t></w:r><w:r><w:t>Cryptography Homework 7 example.docx</w:t></w:r><w:r><w:t>$rq</w:t></w:r></w:p>"
$examplePara.Range.InsertXML($exampleXml)

# ---------------------------------------------------------------------------
# 2) Append a second run to the "Do this lab..." paragraph about working
#    online.
# ---------------------------------------------------------------------------
$doThisIdx2 = Find-ParaIndex("Do this lab in groups of two")
$doThisPara2 = $d.Paragraphs($doThisIdx2)
$doThisXml = "<w:p $wNs><w:r><w:t>Do this lab in groups of two (or three if there is an odd number of students.)</w:t></w:r><w:r><w:t xml:space=`"preserve`">  If you are doing this online, you just have to be both sides of the key exchange.</w:t></w:r></w:p>"
$doThisPara2.Range.InsertXML($doThisXml)

# ---------------------------------------------------------------------------
# 3) Simplify the "Note: ... Andrea Corbillini ..." paragraph: remove the
#    spell-check markup and merge the runs, leaving the hyperlink run itself
#    untouched (so its Hyperlink character-style reference survives).
# ---------------------------------------------------------------------------
$noteIdx = Find-ParaIndex("Note:  The series of four blogs")
$notePara = $d.Paragraphs($noteIdx)
$noteXml = "<w:p $wNs xmlns:r=`"http://schemas.openxmlformats.org/officeDocument/2006/relationships`"><w:r><w:t xml:space=`"preserve`">Note:  The series of four blogs on the subject, </w:t></w:r><w:hyperlink r:id=`"rId5`" w:history=`"1`"><w:r><w:t>Elliptic Curve Cryptography: a gentle introduction</w:t></w:r></w:hyperlink><w:r><w:t xml:space=`"preserve`">, by Andrea Corbillini, is awesome!  We covered some of her basic material from the first blog in class.  The remaining blogs cover finite fields (subgroups and base points, very important), details of ECDH and ECDSA (EC Digital Signature Algorithm), and attacks against discrete logarithm problems.  If you are at all interested</w:t></w:r><w:r><w:t xml:space=`"preserve`"> in ECC</w:t></w:r><w:r><w:t>, her blogs are the place to start.</w:t></w:r></w:p>"
$notePara.Range.InsertXML($noteXml)
# restore the Hyperlink character style that InsertXML cannot set via rStyle
$noteIdx = Find-ParaIndex("Note:  The series of four blogs")
$notePara = $d.Paragraphs($noteIdx)
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("Elliptic Curve Cryptography: a gentle introduction") | Out-Null
$rng.Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 4) Remove the <w:lastRenderedPageBreak/> element from the "The problem of
#    subgroups..." paragraph.
# ---------------------------------------------------------------------------
$subIdx = Find-ParaIndex("The problem of subgroups is common")
$subPara = $d.Paragraphs($subIdx)
$subXml = "<w:p $wNs><w:r><w:t xml:space=`"preserve`">The problem of subgroups is common to all encryption that uses finite fields (i.e., </w:t></w:r><w:r><w:t>most</w:t></w:r><w:r><w:t xml:space=`"preserve`"> of them), not just ECC.</w:t></w:r></w:p>"
$subPara.Range.InsertXML($subXml)

# ---------------------------------------------------------------------------
# 5) Expand "Which point provides the best security?" into several runs.
# ---------------------------------------------------------------------------
$whichIdx = Find-ParaIndex("Which point provides the best security?")
$whichPara = $d.Paragraphs($whichIdx)
$whichXml = "<w:p $wNs><w:pPr><w:pStyle w:val=`"ListParagraph`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr></w:pPr><w:r><w:t>Which point</w:t></w:r><w:r><w:t xml:space=`"preserve`"> on the curve in question 2</w:t></w:r><w:r><w:t xml:space=`"preserve`"> provides the best security</w:t></w:r><w:r><w:t xml:space=`"preserve`">, i.e., has the largest number of points </w:t></w:r><w:r><w:t>?</w:t></w:r></w:p>"
$whichPara.Range.InsertXML($whichXml)

# ---------------------------------------------------------------------------
# 6) Remove the _GoBack bookmark paragraph (two paragraphs after "Which
#    point..."), leaving a plain empty paragraph.
# ---------------------------------------------------------------------------
$whichIdx = Find-ParaIndex("Which point")
$bookmarkPara = $d.Paragraphs($whichIdx + 2)
$bookmarkPara.Range.InsertXML("<w:p $wNs/>")

Write-Output "DONE"
